# Insert a new data row at row 334 (pushing the existing rows 334-426 down
# to 335-427), then populate the new row with its data. This mirrors the
# source diff, which shows a new "Primera" quality record (Fecha 44932 /
# Volumen 1050 / Precio promedio ponderado 676) inserted ahead of the
# existing run of rows for "Vega Monumental Concepción - Betarraga".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 334..426 down to 335..427, duplicating row 334's formatting
# (date style, etc.) onto the freshly inserted row - exactly what Excel does
# for Rows.Insert().
$ws.Rows.Item(334).Insert()

# Populate the newly inserted row 334 with its values.
$ws.Range("A334").Value = 11
$ws.Range("B334").Value = "Vega Monumental Concepción"
$ws.Range("C334").Value = "Bíobío"
$ws.Range("D334").Value = 44932
$ws.Range("E334").Value = 8
$ws.Range("F334").Value = 100114014
$ws.Range("G334").Value = "Betarraga"
$ws.Range("H334").Value = "Sin especificar"
$ws.Range("I334").Value = "Primera"
$ws.Range("J334").Value = 1050
$ws.Range("K334").Value = 650
$ws.Range("L334").Value = 700
$ws.Range("M334").Value = 676
$ws.Range("N334").Value = "`$/paquete 5 unidades"
$ws.Range("O334").Value = "Región Metropolitana"
$ws.Range("P334").Value = 135
$ws.Range("Q334").Value = 5
$ws.Range("R334").Value = "Hortaliza"
